$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the report title (sheet title + workbook "as of" month) and the
#    "Rolling 12 Months" label from October -> November.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Table 1.2.A. Net Generation by Energy Source:  Electric Utilities, 2006-November 2016"

# ---------------------------------------------------------------------------
# 2. Insert a new row for "November" under the "Year 2016" monthly block
#    (it becomes the new row 53, pushing the Year-to-Date / Rolling-12-months
#    blocks down by one row). Applying a full thin border to the inserted
#    row lines it up with the existing "month data" cell style already used
#    by the other month rows (e.g. row 52 / October), instead of minting a
#    brand-new (unbordered) style.
# ---------------------------------------------------------------------------
$ws.Rows.Item(53).Insert()
$ws.Range("A53:M53").Borders.LineStyle = 1

$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 64946
$ws.Range("C53").Value2 = 673
$ws.Range("D53").Value2 = 596
$ws.Range("E53").Value2 = 43773
$ws.Range("F53").Value2 = 22
$ws.Range("G53").Value2 = 33082
$ws.Range("H53").Value2 = 17741
$ws.Range("I53").Value2 = 145
$ws.Range("J53").Value2 = 3631
$ws.Range("K53").Value2 = -522
$ws.Range("L53").Value2 = 24
$ws.Range("M53").Value2 = 164111

# ---------------------------------------------------------------------------
# 3. "Year to Date" block (now rows 55-57) - refresh the 2014/2015/2016
#    totals to include November.
# ---------------------------------------------------------------------------
$ws.Range("B55").Value2 = 1079786
$ws.Range("C55").Value2 = 10029
$ws.Range("D55").Value2 = 8260
$ws.Range("E55").Value2 = 461300
$ws.Range("F55").Value2 = 97
$ws.Range("G55").Value2 = 381444
$ws.Range("H55").Value2 = 217771
$ws.Range("I55").Value2 = 1146
$ws.Range("J55").Value2 = 30529
$ws.Range("K55").Value2 = -4735
$ws.Range("L55").Value2 = 566
$ws.Range("M55").Value2 = 2186194

$ws.Range("B56").Value2 = 929827
$ws.Range("C56").Value2 = 9736
$ws.Range("D56").Value2 = 7674
$ws.Range("E56").Value2 = 565471
$ws.Range("F56").Value2 = 197
$ws.Range("G56").Value2 = 380683
$ws.Range("H56").Value2 = 208344
$ws.Range("I56").Value2 = 1396
$ws.Range("J56").Value2 = 32415
$ws.Range("K56").Value2 = -3895
$ws.Range("L56").Value2 = 509
$ws.Range("M56").Value2 = 2132358

$ws.Range("B57").Value2 = 835856
$ws.Range("C57").Value2 = 7913
$ws.Range("D57").Value2 = 8214
$ws.Range("E57").Value2 = 608453
$ws.Range("F57").Value2 = 142
$ws.Range("G57").Value2 = 387127
$ws.Range("H57").Value2 = 225608
$ws.Range("I57").Value2 = 2033
$ws.Range("J57").Value2 = 36244
$ws.Range("K57").Value2 = -4972
$ws.Range("L57").Value2 = 288
$ws.Range("M57").Value2 = 2106905

# ---------------------------------------------------------------------------
# 4. "Rolling 12 Months Ending in ..." label (now row 58) - October -> November.
# ---------------------------------------------------------------------------
$ws.Range("A58").Value2 = "Rolling 12 Months Ending in November"

# ---------------------------------------------------------------------------
# 5. "Rolling 12 Months" data block (now rows 59-60) - refresh totals.
# ---------------------------------------------------------------------------
$ws.Range("B59").Value2 = 1023113
$ws.Range("C59").Value2 = 10403
$ws.Range("D59").Value2 = 8562
$ws.Range("E59").Value2 = 605586
$ws.Range("F59").Value2 = 212
$ws.Range("G59").Value2 = 419110
$ws.Range("H59").Value2 = 228759
$ws.Range("I59").Value2 = 1468
$ws.Range("J59").Value2 = 35163
$ws.Range("K59").Value2 = -4304
$ws.Range("L59").Value2 = 565
$ws.Range("M59").Value2 = 2328637

$ws.Range("B60").Value2 = 904414
$ws.Range("C60").Value2 = 8562
$ws.Range("D60").Value2 = 8818
$ws.Range("E60").Value2 = 660798
$ws.Range("F60").Value2 = "NM"
$ws.Range("G60").Value2 = 423124
$ws.Range("H60").Value2 = 246904
$ws.Range("I60").Value2 = 2130
$ws.Range("J60").Value2 = 39821
$ws.Range("K60").Value2 = -5182
$ws.Range("L60").Value2 = 336
$ws.Range("M60").Value2 = 2289870

"done"
